$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly fruit/vegetable price data refresh for rows 2-46.
# Each entry: Row, Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L), Precio promedio ponderado (M), Precio $/Kg (P)
$data = @(
    @(2, 44441, 700, 28000, 30000, 29000, 1160),
    @(3, 44504, 600, 11000, 13000, 12000, 480),
    @(4, 44455, 800, 28000, 30000, 29000, 1160),
    @(5, 44419, 600, 27000, 29000, 28000, 1120),
    @(6, 44412, 600, 25000, 27000, 26000, 1040),
    @(7, 44433, 400, 28000, 30000, 29000, 1160),
    @(8, 44462, 400, 22000, 23000, 22500, 900),
    @(9, 44503, 400, 11000, 13000, 12000, 480),
    @(10, 44399, 400, 26000, 28000, 27000, 1080),
    @(11, 44377, 500, 26000, 28000, 27000, 1080),
    @(12, 44497, 500, 13000, 15000, 14000, 560),
    @(13, 44435, 900, 28000, 30000, 29000, 1160),
    @(14, 44356, 300, 26000, 28000, 27000, 1080),
    @(15, 44448, 400, 28000, 30000, 29000, 1160),
    @(16, 44343, 200, 26000, 28000, 27000, 1080),
    @(17, 44363, 240, 28000, 30000, 29000, 1160),
    @(18, 44406, 600, 26000, 28000, 27000, 1080),
    @(19, 44392, 100, 26000, 28000, 27000, 1080),
    @(20, 44426, 400, 28000, 30000, 29000, 1160),
    @(21, 44384, 400, 26000, 28000, 27000, 1080),
    @(22, 44483, 300, 18000, 20000, 19000, 760),
    @(23, 44349, 600, 26000, 28000, 27000, 1080),
    @(24, 44482, 500, 18000, 20000, 19000, 760),
    @(25, 44364, 200, 28000, 30000, 29000, 1160),
    @(26, 44469, 600, 22000, 24000, 23000, 920),
    @(27, 44427, 300, 28000, 30000, 29000, 1160),
    @(28, 44370, 400, 27000, 28000, 27500, 1100),
    @(29, 44475, 1000, 22000, 24000, 23000, 920),
    @(30, 44468, 500, 23000, 25000, 24000, 960),
    @(31, 44489, 400, 18000, 20000, 19000, 760),
    @(32, 44476, 500, 23000, 24000, 23500, 940),
    @(33, 44447, 600, 28000, 30000, 29000, 1160),
    @(34, 44434, 500, 28000, 30000, 29000, 1160),
    @(35, 44490, 500, 16000, 18000, 17000, 680),
    @(36, 44391, 100, 26000, 28000, 27000, 1080),
    @(37, 44420, 700, 27000, 29000, 28000, 1120),
    @(38, 44385, 500, 26000, 28000, 27000, 1080),
    @(39, 44350, 700, 28000, 30000, 29000, 1160),
    @(40, 44405, 500, 26000, 28000, 27000, 1080),
    @(41, 44413, 700, 26000, 28000, 27000, 1080),
    @(42, 44398, 500, 26000, 28000, 27000, 1080),
    @(43, 44461, 500, 23000, 25000, 24000, 960),
    @(44, 44357, 340, 28000, 30000, 29000, 1160),
    @(45, 44371, 500, 28000, 30000, 29000, 1160),
    @(46, 44454, 1000, 28000, 30000, 29000, 1160)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $row[2]   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $row[3]   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[4]   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[5]   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row[6]   # P: Precio $/Kg
}

